$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 714 (shifts existing rows 714:736 down to 715:737),
# carrying a new weekly price record for "Primera" quality pineapple.
$ws.Rows("714:714").Insert()

$ws.Cells.Item(714, 1).Value = 10
$ws.Cells.Item(714, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(714, 3).Value = "La Araucanía"
$ws.Cells.Item(714, 4).Value = 45075
$ws.Cells.Item(714, 5).Value = 9
$ws.Cells.Item(714, 6).Value = "Fruta"
$ws.Cells.Item(714, 7).Value = 100108
$ws.Cells.Item(714, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(714, 9).Value = 100108005
$ws.Cells.Item(714, 10).Value = "Piña"
$ws.Cells.Item(714, 11).Value = "Caramelo"
$ws.Cells.Item(714, 12).Value = "Primera"
$ws.Cells.Item(714, 13).Value = 180
$ws.Cells.Item(714, 14).Value = 20000
$ws.Cells.Item(714, 15).Value = 20000
$ws.Cells.Item(714, 16).Value = 20000
$ws.Cells.Item(714, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(714, 18).Value = "Ecuador"
$ws.Cells.Item(714, 19).Value = 1667
$ws.Cells.Item(714, 20).Value = 12
